# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - text first, then copy the existing header's formatting
# (bold / bordered / centered style already used by the rest of row 1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every player row (2 through 45)
$wins = 86
$losses = 76
$ties = 0

for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD = 30
    $ws.Cells.Item($row, 31).Value = $losses  # column AE = 31
    $ws.Cells.Item($row, 32).Value = $ties    # column AF = 32
}
